$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.887.11'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.623.99'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.39'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.516'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.05'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.34%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0881'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.853.74'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.632.71'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.553'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.52'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.890.55'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '227.38'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.93'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.25'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.91'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.37'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.417.52'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.985'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '65.30'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.37'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.79'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.763.69'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '89.26'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.45%  '
